# Apply updates described by the diff to sheet1 (rows 12-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 ---
$ws.Range("B12").Value = 89950

# --- Row 13 ---
$ws.Range("A13").Value = 112128708
$ws.Range("B13").Value = 90835
$ws.Range("E13").Value = 5964
$ws.Range("F13").Value = "Fjällig taggsvamp s.str."
$ws.Range("G13").Value = "Sarcodon imbricatus s.str."
$ws.Range("H13").Value = "(L.:Fr.) P.Karst."
# I13 holds a numeric-looking count ("1") that must stay text, not be
# coerced to a number - force text formatting for the assignment, then
# drop back to the Normal style so no stray formatting is left behind.
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "1"
$ws.Range("I13").Style = "Normal"
$ws.Range("J13").Value = "fruktkroppar"
$ws.Range("Q13").Value = 657216
$ws.Range("R13").Value = 6571313

# --- Row 14 ---
$ws.Range("A14").Value = 112128627
$ws.Range("B14").Value = 90235
$ws.Range("E14").Value = 3298
$ws.Range("F14").Value = "Trådticka"
$ws.Range("G14").Value = "Climacocystis borealis"
$ws.Range("H14").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I14").ClearContents()
$ws.Range("J14").ClearContents()
$ws.Range("Q14").Value = 657182
$ws.Range("R14").Value = 6571192
$ws.Range("AC14").Value = "På nedre delen av torrgran."

# --- Row 15 ---
$ws.Range("A15").Value = 112128573
$ws.Range("B15").Value = 90816
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 788
$ws.Range("F15").Value = "Gul taggsvamp"
$ws.Range("G15").Value = "Hydnellum geogenium"
$ws.Range("H15").Value = "(Fr.) Banker"
$ws.Range("Q15").Value = 657134
$ws.Range("R15").Value = 6571219
$ws.Range("AC15").ClearContents()

# --- Row 16 ---
$ws.Range("A16").Value = 112128498
$ws.Range("B16").Value = 90166
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 1339
$ws.Range("F16").Value = "Brandticka"
$ws.Range("G16").Value = "Pycnoporellus fulgens"
$ws.Range("H16").Value = "(Fr.) Donk"
$ws.Range("R16").Value = 6571271
$ws.Range("AC16").Value = "På granlåga. En del årsfärska dödade granar av granbarkborre. Gott om död ved i form av torrträd och lågor av gran."

# --- Row 17 ---
$ws.Range("A17").Value = 112128551
$ws.Range("B17").Value = 90814
$ws.Range("E17").Value = 4364
$ws.Range("F17").Value = "Dropptaggsvamp"
$ws.Range("G17").Value = "Hydnellum ferrugineum"
$ws.Range("H17").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q17").Value = 657162
$ws.Range("AC17").ClearContents()

Write-Host "Applied updates to rows 12-17"
